$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the EVENTO number in A2
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1553046

# ---------------------------------------------------------------------------
# 2) Remove the per-column direct formatting ("style" attribute) that used
#    to be applied to columns B, F and I, while preserving the formatting
#    that is actually applied on the individual cells that already contain
#    data (header row 1, and the date values in column B).
#    We do this by stashing the current cell-level formats of representative
#    cells into scratch cells far outside the used range, clearing the
#    column-level formatting, and then pasting the stashed formats back onto
#    the cells that need to keep them.
# ---------------------------------------------------------------------------
$ws.Range("D1").Copy()                  # D1 carries the same header style as B1/F1/I1 (untouched column)
$ws.Range("Z1").PasteSpecial(-4122)     # -4122 = xlPasteFormats

$ws.Range("B2").Copy()                  # B2 carries the date style used across B2:B8
$ws.Range("Z2").PasteSpecial(-4122)

$ws.Columns("B:B").ClearFormats()
$ws.Columns("F:F").ClearFormats()
$ws.Columns("I:I").ClearFormats()

$ws.Range("Z1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("Z2").Copy()
$ws.Range("B2:B8").PasteSpecial(-4122)

$ws.Range("Z1:Z2").Clear()

# ---------------------------------------------------------------------------
# 3) Clear the event rows 3 to 8 (only column B keeps its empty, date
#    formatted cell; every other column becomes fully empty again).
# ---------------------------------------------------------------------------
$ws.Range("A3:I8").ClearContents()

# ---------------------------------------------------------------------------
# 4) Update the selection to B3:I8 with active cell B3 (mimics the user
#    entering the EVENTO screen and selecting the editable area).
# ---------------------------------------------------------------------------
$ws.Range("B3:I8").Select() | Out-Null
